$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume figures (refreshed data snapshot).
# Column D prices are stored as text (e.g. "1.005", "28.718.39") rather than
# numbers, so we force a text number format on each price cell before writing
# its value - this stops Excel from reinterpreting values like "1.005" as the
# number 1.005 (which would drop the trailing zero) or similar coercions.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.718.39"
$ws.Range("E2").Value = "  +2.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.874.42"
$ws.Range("E3").Value = "  +2.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.27"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4595"
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07861"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9978"
$ws.Range("E10").Value = "  +3.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.77"
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.893.75"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.985"
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.702"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06955"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.42"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.85"
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.729.91"
$ws.Range("E21").Value = "  +2.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.276"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.02"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.124"
$ws.Range("E24").Value = "  +1.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.098.07"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.73"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.19"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.767"
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.961"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.00"
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09321"
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9193"
$ws.Range("E32").Value = "  -2.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.306"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.320"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05755"
$ws.Range("E36").Value = "  -1.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.148"
$ws.Range("E37").Value = "  +1.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02074"
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.719"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5635"
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1788"
$ws.Range("E41").Value = "  +1.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.906"
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07206"
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.69"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5282"
$ws.Range("E45").Value = "  +0.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.136"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "113.65"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.826"
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.408"
$ws.Range("E50").Value = "  +4.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.004"
$ws.Range("E51").Value = "  +0.36%  "
